{"js": "// Replace the date line and every arithmetic-expression cell in the\n// practice-sheet table with the new values from the target revision.\n// The document's paragraph order is: [0] the title/date paragraph,\n// followed by [1..100] one paragraph per table cell (row-major, 5 cols\n// x 20 rows) -- this matches the order the source diff lists them in.\nconst oldTexts = [\"2023-09-27 Wednesday\", \"50+12=\", \"92-60=\", \"9+12=\", \"98-14=\", \"78-34=\", \"37+42=\", \"15+82=\", \"72-37=\", \"52-28=\", \"53+26=\", \"11+9=\", \"18+19=\", \"69-62=\", \"37+4=\", \"22+48=\", \"35+33=\", \"3+38=\", \"93-36=\", \"10+83=\", \"16-11=\", \"34+62=\", \"89-77=\", \"96-67=\", \"51-27=\", \"46-9=\", \"44-16=\", \"59-3=\", \"59-39=\", \"31-20=\", \"74-59=\", \"16+8=\", \"78-9=\", \"56-42=\", \"35+7=\", \"16+46=\", \"65+4=\", \"59+20=\", \"26+44=\", \"71-1=\", \"7+23=\", \"95-89=\", \"79-5=\", \"60+25=\", \"50-23=\", \"30+69=\", \"8+48=\", \"87-58=\", \"25+27=\", \"73-46=\", \"45+49=\", \"96-72=\", \"72+8=\", \"74+10=\", \"88-9=\", \"86-54=\", \"69+23=\", \"11+59=\", \"37+19=\", \"6+53=\", \"69-69=\", \"93-5=\", \"56-52=\", \"50-3=\", \"54-32=\", \"74+3=\", \"90-56=\", \"27+26=\", \"55+35=\", \"45-18=\", \"36+44=\", \"74+25=\", \"80-50=\", \"38-22=\", \"41-14=\", \"27+59=\", \"95-72=\", \"11+22=\", \"49-4=\", \"31+20=\", \"54+38=\", \"95-4=\", \"73-25=\", \"57-42=\", \"62-10=\", \"27+13=\", \"62-33=\", \"57+40=\", \"65-30=\", \"50-30=\", \"17+78=\", \"78-66=\", \"10-5=\", \"12-1=\", \"0+97=\", \"46-3=\", \"54+39=\", \"20+1=\", \"16+24=\", \"6+13=\", \"21-6=\"];\nconst newTexts = [\"2023-09-28 Thursday\", \"62-32=\", \"95-42=\", \"76-57=\", \"89-5=\", \"99-25=\", \"98-88=\", \"85-71=\", \"94-2=\", \"14+1=\", \"3-2=\", \"20+25=\", \"42+31=\", \"78-32=\", \"32+58=\", \"12+18=\", \"24+52=\", \"36-12=\", \"64-53=\", \"58-45=\", \"29+29=\", \"94-85=\", \"72+6=\", \"60-2=\", \"85-50=\", \"91-61=\", \"94-12=\", \"53+5=\", \"64-37=\", \"21+55=\", \"95-40=\", \"42+30=\", \"6+65=\", \"38+22=\", \"75-37=\", \"43+3=\", \"7-2=\", \"61-20=\", \"8+76=\", \"42-13=\", \"58+34=\", \"55+21=\", \"21+57=\", \"29+55=\", \"63+19=\", \"82-66=\", \"13+15=\", \"31-22=\", \"27+58=\", \"9+4=\", \"62-51=\", \"83-71=\", \"84-30=\", \"77-0=\", \"95-25=\", \"25-4=\", \"51+21=\", \"74+24=\", \"34+43=\", \"96-85=\", \"59+21=\", \"60-28=\", \"8-6=\", \"80-58=\", \"13+35=\", \"54+23=\", \"41+48=\", \"24+19=\", \"30+13=\", \"46+26=\", \"52-9=\", \"93-45=\", \"9+61=\", \"2+95=\", \"99-87=\", \"4+4=\", \"35-23=\", \"30+41=\", \"85-84=\", \"96-20=\", \"62-54=\", \"18+14=\", \"59+34=\", \"98-37=\", \"56-35=\", \"34-26=\", \"23+2=\", \"66+6=\", \"58-20=\", \"50-8=\", \"4-2=\", \"23-12=\", \"79-36=\", \"94-12=\", \"16-7=\", \"82-70=\", \"6+56=\", \"97-70=\", \"90+9=\", \"70-52=\", \"54-46=\"];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst count = Math.min(paragraphs.items.length, newTexts.length);\nfor (let i = 0; i < count; i++) {\n  const para = paragraphs.items[i];\n  const current = para.text;\n  const expectedOld = oldTexts[i];\n  const newText = newTexts[i];\n  // Only touch paragraphs whose text actually changes, and guard against\n  // drift between our recorded order and the live document by checking\n  // the paragraph still holds the value we expect to replace.\n  if (current === newText) {\n    continue;\n  }\n  if (current === expectedOld) {\n    para.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and every arithmetic-expression cell in the\n# practice-sheet table to the new values from the target revision.\n# Cell values are listed row-major (row 1 col 1..5, row 2 col 1..5, ...)\n# matching Table.Cell($row, $col) iteration order.\n\n$d = $word.ActiveDocument\n\n# 1. Title/date line (plain paragraph above the table).\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"2023-09-27 Wednesday\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"2023-09-28 Thursday\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# 2. Table cells.\n$oldValues = @(\"50+12=\", \"92-60=\", \"9+12=\", \"98-14=\", \"78-34=\", \"37+42=\", \"15+82=\", \"72-37=\", \"52-28=\", \"53+26=\", \"11+9=\", \"18+19=\", \"69-62=\", \"37+4=\", \"22+48=\", \"35+33=\", \"3+38=\", \"93-36=\", \"10+83=\", \"16-11=\", \"34+62=\", \"89-77=\", \"96-67=\", \"51-27=\", \"46-9=\", \"44-16=\", \"59-3=\", \"59-39=\", \"31-20=\", \"74-59=\", \"16+8=\", \"78-9=\", \"56-42=\", \"35+7=\", \"16+46=\", \"65+4=\", \"59+20=\", \"26+44=\", \"71-1=\", \"7+23=\", \"95-89=\", \"79-5=\", \"60+25=\", \"50-23=\", \"30+69=\", \"8+48=\", \"87-58=\", \"25+27=\", \"73-46=\", \"45+49=\", \"96-72=\", \"72+8=\", \"74+10=\", \"88-9=\", \"86-54=\", \"69+23=\", \"11+59=\", \"37+19=\", \"6+53=\", \"69-69=\", \"93-5=\", \"56-52=\", \"50-3=\", \"54-32=\", \"74+3=\", \"90-56=\", \"27+26=\", \"55+35=\", \"45-18=\", \"36+44=\", \"74+25=\", \"80-50=\", \"38-22=\", \"41-14=\", \"27+59=\", \"95-72=\", \"11+22=\", \"49-4=\", \"31+20=\", \"54+38=\", \"95-4=\", \"73-25=\", \"57-42=\", \"62-10=\", \"27+13=\", \"62-33=\", \"57+40=\", \"65-30=\", \"50-30=\", \"17+78=\", \"78-66=\", \"10-5=\", \"12-1=\", \"0+97=\", \"46-3=\", \"54+39=\", \"20+1=\", \"16+24=\", \"6+13=\", \"21-6=\")\n$newValues = @(\"62-32=\", \"95-42=\", \"76-57=\", \"89-5=\", \"99-25=\", \"98-88=\", \"85-71=\", \"94-2=\", \"14+1=\", \"3-2=\", \"20+25=\", \"42+31=\", \"78-32=\", \"32+58=\", \"12+18=\", \"24+52=\", \"36-12=\", \"64-53=\", \"58-45=\", \"29+29=\", \"94-85=\", \"72+6=\", \"60-2=\", \"85-50=\", \"91-61=\", \"94-12=\", \"53+5=\", \"64-37=\", \"21+55=\", \"95-40=\", \"42+30=\", \"6+65=\", \"38+22=\", \"75-37=\", \"43+3=\", \"7-2=\", \"61-20=\", \"8+76=\", \"42-13=\", \"58+34=\", \"55+21=\", \"21+57=\", \"29+55=\", \"63+19=\", \"82-66=\", \"13+15=\", \"31-22=\", \"27+58=\", \"9+4=\", \"62-51=\", \"83-71=\", \"84-30=\", \"77-0=\", \"95-25=\", \"25-4=\", \"51+21=\", \"74+24=\", \"34+43=\", \"96-85=\", \"59+21=\", \"60-28=\", \"8-6=\", \"80-58=\", \"13+35=\", \"54+23=\", \"41+48=\", \"24+19=\", \"30+13=\", \"46+26=\", \"52-9=\", \"93-45=\", \"9+61=\", \"2+95=\", \"99-87=\", \"4+4=\", \"35-23=\", \"30+41=\", \"85-84=\", \"96-20=\", \"62-54=\", \"18+14=\", \"59+34=\", \"98-37=\", \"56-35=\", \"34-26=\", \"23+2=\", \"66+6=\", \"58-20=\", \"50-8=\", \"4-2=\", \"23-12=\", \"79-36=\", \"94-12=\", \"16-7=\", \"82-70=\", \"6+56=\", \"97-70=\", \"90+9=\", \"70-52=\", \"54-46=\")\n\n$t = $d.Tables.Item(1)\n$i = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $expectedOld = $oldValues[$i]\n        $newValue = $newValues[$i]\n        $cellRange = $cell.Range\n        $cellRange.MoveEnd(1, -1) | Out-Null\n        if ($cellRange.Text -eq $expectedOld) {\n            $cellRange.Text = $newValue\n        }\n        $i++\n    }\n}\n"}
